$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31: num_customers 42 -> 45, retention_rate recalculated (45/2312)
$ws.Range("C31").Value = 45
$ws.Range("E31").Value = 45/2312

# Row 36: num_customers 118 -> 119, retention_rate recalculated (119/1930)
$ws.Range("C36").Value = 119
$ws.Range("E36").Value = 119/1930

# Row 37: num_customers 725 -> 731, cohort_size 725 -> 731, retention_rate stays 1
$ws.Range("C37").Value = 731
$ws.Range("D37").Value = 731
